$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are digit-only strings (e.g. "1.00", "559.22") that must
# stay plain text (as they were stored originally) instead of being
# auto-converted to numbers by the input parser. Temporarily mark the cell
# as Text before assigning, then clear the format again so the cell keeps
# the default (General) style, matching the original file's formatting.

$ws.Range("D2").Value = '65.795.00'
$ws.Range("E2").Value = '  -5.57%  '
$ws.Range("D3").Value = '3.277.57'
$ws.Range("E3").Value = '  -6.12%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.89'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.19%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.22%  '
$ws.Range("D9").Value = '3.271.65'
$ws.Range("E9").Value = '  -5.91%  '
$ws.Range("E10").Value = '  -9.70%  '
$ws.Range("E11").Value = '  -5.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.36'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -8.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000266'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -7.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '638.25'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.62'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.76%  '
$ws.Range("D16").Value = '3.798.64'
$ws.Range("E16").Value = '  -5.66%  '
$ws.Range("D17").Value = '65.781.42'
$ws.Range("E17").Value = '  -5.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.96'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.54%  '
$ws.Range("E19").Value = '  -3.36%  '
$ws.Range("D20").Value = '3.283.74'
$ws.Range("E20").Value = '  -5.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.33'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -8.11%  '
$ws.Range("E22").Value = '  -4.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.26'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '107.37'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +8.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.90'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -6.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.97'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -7.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.68'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -7.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.52'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -5.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.69'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -7.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.31'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -7.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.98'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.25'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.02'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.30%  '
$ws.Range("E34").Value = '  -4.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.67'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.60%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '3.720.69'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '523.35'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.38'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.49%  '
$ws.Range("D40").Value = '0.0₃0730'
$ws.Range("E40").Value = '  -7.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.130'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '32.95'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.36'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.53%  '
$ws.Range("E45").Value = '  -10.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.27'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("E47").Value = '  -6.72%  '
$ws.Range("E48").Value = '  -4.15%  '
$ws.Range("E49").Value = '  -8.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("E51").Value = '  +2.38%  '
